$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$participantsQuery = @'
with diagnosis1 as (
select dgn."participant.id", group_concat(dgn.age_at_diagnosis,';') as age, group_concat(dgn.diagnosis,';') as diag,group_concat(dgn.anatomic_site,';') as ant_site from df_diagnosis dgn where dgn."participant.id" is not null group by dgn."participant.id" ),
diagnosis2 as (select "participant.id",  group_concat(diagnosis,';') as diag from (select distinct "participant.id", diagnosis from df_diagnosis  where "participant.id" is not null )  group by "participant.id" ),
diagnosis3 as (select "participant.id",  group_concat(anatomic_site,';') as ant_site from (select distinct "participant.id", anatomic_site from df_diagnosis where "participant.id" is not null ) group by "participant.id" ),
treatment1 as (select trt."participant.id",trt.treatment_type from  df_treatment trt  where trt.treatment_type is not null)
SELECT DISTINCT
    prt.participant_id AS "Participant ID",
    std.dbgap_accession AS "Study ID",
    COALESCE(prt.sex_at_birth, '') AS "Sex",
    COALESCE(prt.race, '') AS "Race",
	dgn2.diag AS "Diagnosis",
	dgn3.ant_site AS "Diagnosis Anatomic Site",
	   COALESCE(CASE WHEN dgn1.age = '-999' THEN 'Not Reported' ELSE dgn1.age END, "") AS "Age at Diagnosis (days)",
	trt1.treatment_type AS "Treatment Type",
	srv.last_known_survival_status AS "Last Known Survival Status"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    diagnosis1 dgn1 ON prt.id = dgn1."participant.id" 
LEFT JOIN 
    diagnosis2 dgn2 ON prt.id = dgn2."participant.id"
LEFT JOIN 
    diagnosis3 dgn3 ON prt.id = dgn3."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_diagnosis dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    treatment1 trt1 ON prt.id = trt1."participant.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
WHERE 
    std.dbgap_accession = 'phs002371' and prt.sex_at_birth='Female' and dgn.diagnosis='9861/3 : Acute myeloid leukemia, NOS' and seq.file_type ='fastq'
ORDER BY 
    prt.participant_id ASC 
'@
$samplesQuery = @'
WITH filtered_participants AS (
    SELECT DISTINCT prt.id, prt.participant_id
    FROM df_study std
    JOIN df_participant prt ON std.id = prt."study.id"
    JOIN df_diagnosis dgn ON prt.id = dgn."participant.id"
    WHERE 
        std.dbgap_accession = 'phs002371'
        AND prt.sex_at_birth = 'Female'
        AND dgn.diagnosis = '9861/3 : Acute myeloid leukemia, NOS'
),
-- Step 2: Get samples from those participants
filtered_samples AS (
    SELECT DISTINCT smp.id, smp.sample_id, smp.anatomic_site, smp."participant.id"
    FROM df_sample smp
    JOIN filtered_participants prt ON smp."participant.id" = prt.id
),
-- Step 3: Keep only samples that have at least one FASTQ file
samples_with_fastq AS (
    SELECT DISTINCT smp.*
    FROM filtered_samples smp
    JOIN df_sequencing_file seq ON seq."sample.id" = smp.id
    WHERE seq.file_type = 'fastq'
)
-- Final output
SELECT 
    smp.sample_id AS "Sample ID",
    prt.participant_id AS "Participant ID",
    'phs002371' AS "Study ID",
    smp.anatomic_site AS "Sample Anatomic Site",
    NULL AS "Sample Diagnosis"
FROM 
    samples_with_fastq smp
JOIN filtered_participants prt ON smp."participant.id" = prt.id
ORDER BY 
    smp.sample_id ASC;
'@
$statQuery = @'
WITH filtered_participants AS (
    SELECT DISTINCT prt.id
    FROM df_study std
    JOIN df_participant prt ON std.id = prt."study.id"
    JOIN df_diagnosis dgn ON prt.id = dgn."participant.id"
    WHERE 
        std.dbgap_accession = 'phs002371'
        AND prt.sex_at_birth = 'Female'
        AND dgn.diagnosis = '9861/3 : Acute myeloid leukemia, NOS'
),
filtered_samples AS (
    SELECT DISTINCT smp.id
    FROM df_sample smp
    JOIN filtered_participants prt ON smp."participant.id" = prt.id
),
filtered_files AS (
    SELECT DISTINCT seq.id
    FROM df_sequencing_file seq
    JOIN filtered_samples smp ON seq."sample.id" = smp.id
    WHERE seq.file_type = 'fastq'
)
SELECT 
    1 AS "Studies",
    (SELECT COUNT(*) FROM filtered_participants) AS "Participants",
    (SELECT COUNT(*) FROM filtered_samples) AS "Samples",
    (SELECT COUNT(*) FROM filtered_files) AS "Files";
'@
$studiesQuery = @'
WITH Study AS (
    SELECT DISTINCT
        std.study_name AS "Study_Name",
        std.dbgap_accession AS "Study_ID",
        std.study_status AS "Study_Status"
    FROM df_study std 
    WHERE std.dbgap_accession IS NOT NULL
),
diagnosis AS (
    SELECT study_id, group_concat(diagcounts, ' ') AS diagcount 
    FROM (
        SELECT 'phs002371' AS study_id, diagnosis || '(' || d_count || ')' AS diagcounts 
        FROM (
            SELECT row_number() OVER() AS row_num, dgn.diagnosis, COUNT(*) AS d_count 
            FROM df_diagnosis dgn 
            GROUP BY dgn.diagnosis 
            ORDER BY d_count DESC
        ) 
        WHERE row_num <= 5
    )
),
diagnosisas AS (
    SELECT study_id, group_concat(anatcounts, ' ') AS anatcount 
    FROM (
        SELECT 'phs002371' AS study_id, anat || '(' || as_count || ')' AS anatcounts 
        FROM (
            SELECT row_number() OVER() AS row_num, dgn1.anatomic_site AS anat, COUNT(*) AS as_count 
            FROM df_diagnosis dgn1 
            GROUP BY dgn1.anatomic_site 
            ORDER BY as_count DESC
        ) 
        WHERE row_num <= 5
    )
),
participants AS (
    SELECT 'phs002371' AS study_id, COUNT(*) AS p_count 
    FROM df_participant prt
),
sample AS (
    SELECT 'phs002371' AS study_id, COUNT(*) AS s_count 
    FROM df_sample
),
total_files AS (
    SELECT 'phs002371' AS study_id, COUNT(*) AS file_count 
    FROM df_sequencing_file
),
files AS (
    SELECT study_id, group_concat(filecounts, ' ') AS file 
    FROM (
        SELECT 'phs002371' AS study_id, file_type || '(' || filecount || ')' AS filecounts 
        FROM (
            SELECT row_number() OVER() AS row_num, file_type, COUNT(*) AS filecount 
            FROM df_sequencing_file 
            GROUP BY file_type 
            ORDER BY filecount DESC
        ) 
        WHERE row_num <= 5
    )
),
publications AS (
    SELECT 'phs002371' AS study_id, group_concat(pubmed_id, ';') AS pubmeds 
    FROM df_publication
),
personnel AS (
    SELECT 'phs002371' AS study_id, group_concat(personnel_name, ';') AS name 
    FROM df_study_personnel
),
funding AS (
    SELECT 'phs002371' AS study, group_concat(grant_id, ';') AS grant 
    FROM (
        SELECT DISTINCT grant_id 
        FROM df_study_funding
    )
)
SELECT 
    s.study_name AS "Study Name",
    s.study_id AS "Study ID",
    s.study_status AS "Study Status",
    d.diagcount AS "Diagnosis (Top 5)",
    da.anatcount AS "Diagnosis Anatomic Site (Top 5)", 
    p.p_count AS "Number of Participants",
    smp.s_count AS "Number of Samples",
    tf.file_count AS "Number of Files",
    f.file AS "File Type (Top 5)",
    pub.pubmeds AS "PubMed ID",
    pr.name AS "Principal Investigator(s)",
    fund.grant AS "Grant ID"
FROM Study s
LEFT JOIN diagnosis d ON s.study_id = d.study_id
LEFT JOIN diagnosisas da ON s.study_id = da.study_id
LEFT JOIN participants p ON s.study_id = p.study_id
LEFT JOIN sample smp ON s.study_id = smp.study_id
LEFT JOIN total_files tf ON s.study_id = tf.study_id
LEFT JOIN files f ON s.study_id = f.study_id
LEFT JOIN publications pub ON s.study_id = pub.study_id
LEFT JOIN personnel pr ON s.study_id = pr.study_id
LEFT JOIN funding fund ON s.study_id = fund.study;
'@

$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("B3").Value = $studiesQuery
$ws.Range("B4").Value = $samplesQuery
